# "modify W10 and add W11"
# W10 salary/task sheet update:
#   - B3 (the reporting period "Date" cell) switches from a literal date
#     value to a free-text week-range label.
#   - B20 gains the next "task to complete" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "2019.11.14 - 2019.11.21"
$ws.Range("B20").Value = "Present the User Study to stakeholders"

# Restore the cursor/selection position to where the author left it.
$ws.Range("B20").Select()
